$d = $word.ActiveDocument

# --- Update the contact/header block ---
$d.Content.Find.Execute("Valencia, Spain", $true, $false, $false, $false, $false, $true, 1, $false, "Valencia, US", 2)
$d.Content.Find.Execute("guillermo@test.com", $true, $false, $false, $false, $false, $true, 1, $false, "guillermo@g.com", 2)
$d.Content.Find.Execute("248-111-2222", $true, $false, $false, $false, $false, $true, 1, $false, "242555111", 2)

# --- Rewrite the cover letter body paragraphs ---

# 1st body paragraph
$d.Content.Find.Execute("I am excited to apply for the Junior Software Engineer position at your company. With my experience as a Software Engineer Intern at Masetto Logistics and Barracuda, I have honed my skills in Python, React, Java, and JavaScript, aligning perfectly with the qualifications you are seeking. I have a strong foundation in database systems, data structures, and algorithms, and I am well-versed in object-oriented design principles.", $true, $false, $false, $false, $false, $true, 1, $false, "I am writing to express my interest in the Windows Engineer position at Epic as advertised. With a background in software engineering and a strong focus on backend technologies, I am confident that my skills and experience align well with the requirements of the role.", 2)

# 2nd body paragraph
$d.Content.Find.Execute("During my internship at Masetto Logistics, I developed real-time fleet management features using Python and React, incorporating API calls for precise location tracking. I collaborated closely with cross-functional teams in agile environments, similar to the collaborative atmosphere described in your role description. My experiences have equipped me with strong problem-solving abilities and the capacity to work both independently and as part of a team.", $true, $false, $false, $false, $false, $true, 1, $false, "My internship experiences in software engineering have equipped me with expertise in automation using scripting tools like Python and Selenium, which I believe can be directly applied to the automation requirements of the position. Additionally, my understanding of configuration management tools such as Puppet and Chef, as well as my familiarity with REST API integration, make me well-suited for the responsibilities of designing and implementing domain architecture, including Active Directory and DNS.", 2)

# 3rd body paragraph - replaced and split into two paragraphs of text (separated by
# a blank line, i.e. two manual line breaks) using the "^l" replace special-character,
# which Word's Find/Replace expands into extra <w:br/> + <w:t> runs.
$d.Content.Find.Execute("I am eager to bring my technical skills and passion for software development to your team. I am confident that my background aligns well with the requirements of the Junior Software Engineer role, and I am excited about the opportunity to contribute to your organization's success. Thank you for considering my application.", $true, $false, $false, $false, $false, $true, 1, $false, "I am particularly excited about the opportunity to work collaboratively with experts in different knowledge areas at Epic to ensure a robust and secure hosting environment. My agile mindset, acquired through experience in sprint planning and scrum meetings, will enable me to adapt quickly to the dynamic challenges of the role. I am prepared to relocate to the Madison, WI area and eager to contribute to the continuous learning environment at Epic.^l^lThank you for considering my application. I am looking forward to the possibility of contributing my skills to the impactful work at Epic. ", 2)
